$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fleet_definition")

# --- New shared string "#" used by column J in every added row ---
# (created implicitly the first time it is assigned below)

# --- Rows 34-52: fleet rows, columns F-L, all cells default style ---
$ws.Range("F34").Value = 24
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = "#"
$ws.Range("K34").Value = 1
$ws.Range("L34").Value = "F1_JPN_WCNPO_OSDWLL_early_Area1  # 1"

$ws.Range("F35").Value = 24
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = "#"
$ws.Range("K35").Value = 2
$ws.Range("L35").Value = "F2_JPN_WCNPO_OSDWCOLL_late_Area1  # 2"

$ws.Range("F36").Value = 24
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = "#"
$ws.Range("K36").Value = 3
$ws.Range("L36").Value = "F3_JPN_EPO_OSDWLL #3"

$ws.Range("F37").Value = 15
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 3
$ws.Range("J37").Value = "#"
$ws.Range("K37").Value = 4
$ws.Range("L37").Value = "F4_JPN_WCNPO_OSDF  # 4"

$ws.Range("F38").Value = 24
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = "#"
$ws.Range("K38").Value = 5
$ws.Range("L38").Value = "F5_JPN_WCNPO_CODF  # 5"

$ws.Range("F39").Value = 15
$ws.Range("G39").Value = 0
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 1
$ws.Range("J39").Value = "#"
$ws.Range("K39").Value = 6
$ws.Range("L39").Value = "F6_JPN_WCNPO_Other_early  # 6"

$ws.Range("F40").Value = 15
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 2
$ws.Range("J40").Value = "#"
$ws.Range("K40").Value = 7
$ws.Range("L40").Value = "F7_JPN_WCNPO_Other_late  # 7"

$ws.Range("F41").Value = 24
$ws.Range("G41").Value = 0
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = "#"
$ws.Range("K41").Value = 8
$ws.Range("L41").Value = "F8_TWN_WCNPO_DWLL_late  # 8"

$ws.Range("F42").Value = 15
$ws.Range("G42").Value = 0
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 8
$ws.Range("J42").Value = "#"
$ws.Range("K42").Value = 9
$ws.Range("L42").Value = "F9_TWN_WCNPO_DWLL_early  # 9"

$ws.Range("F43").Value = 15
$ws.Range("G43").Value = 0
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 2
$ws.Range("J43").Value = "#"
$ws.Range("K43").Value = 10
$ws.Range("L43").Value = "F10_TWN_WCNPO_Other  # 10"

$ws.Range("F44").Value = 24
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = "#"
$ws.Range("K44").Value = 11
$ws.Range("L44").Value = "F11_US_WCNPO_LL_deep  # 11"

$ws.Range("F45").Value = 24
$ws.Range("G45").Value = 0
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = "#"
$ws.Range("K45").Value = 12
$ws.Range("L45").Value = "F12_US_WCNPO_LL_shallow_late  # 12"

$ws.Range("F46").Value = 24
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = "#"
$ws.Range("K46").Value = 13
$ws.Range("L46").Value = "F13_US_WCNPO_LL_shallow_early  # 13"

$ws.Range("F47").Value = 15
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 8
$ws.Range("J47").Value = "#"
$ws.Range("K47").Value = 14
$ws.Range("L47").Value = "F14_US_WCNPO_GN  # 14"

$ws.Range("F48").Value = 15
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 8
$ws.Range("J48").Value = "#"
$ws.Range("K48").Value = 15
$ws.Range("L48").Value = "F15_US_WCNPO_Other  # 15"

$ws.Range("F49").Value = 15
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 11
$ws.Range("J49").Value = "#"
$ws.Range("K49").Value = 16
$ws.Range("L49").Value = "F16_JPN_WCNPO_OSDWLL_early_Area2  # 16"

$ws.Range("F50").Value = 15
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 11
$ws.Range("J50").Value = "#"
$ws.Range("K50").Value = 17
$ws.Range("L50").Value = "F17_JPN_WCNPO_OSDWLL_late_Area2  # 17"

$ws.Range("F51").Value = 15
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 8
$ws.Range("J51").Value = "#"
$ws.Range("K51").Value = 18
$ws.Range("L51").Value = "F18_WCPFC # 18"

$ws.Range("F52").Value = 24
$ws.Range("G52").Value = 0
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = "#"
$ws.Range("K52").Value = 19
$ws.Range("L52").Value = "F19_IATTC  # 19"

# --- Rows 53-60: sensitivity rows, columns F-L (no I), L uses the "Fleet name"-like bold/explicit font style ---
$ws.Range("F53").Value = 15
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = "#"
$ws.Range("K53").Value = 20
$ws.Range("L53").Value = "S1_JPN_WCNPO_OSDWLL_early_Area1"
$ws.Range("L53").Font.ThemeColor = 1

$ws.Range("F54").Value = 15
$ws.Range("G54").Value = 0
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = "#"
$ws.Range("K54").Value = 21
$ws.Range("L54").Value = "S2_JPN_WCNPO_OSDWCOLL_late_Area1"
$ws.Range("L54").Font.ThemeColor = 1

$ws.Range("F55").Value = 15
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = "#"
$ws.Range("K55").Value = 22
$ws.Range("L55").Value = "S3_JPN_WCNPO_OSDWLL_early_Area2"
$ws.Range("L55").Font.ThemeColor = 1

$ws.Range("F56").Value = 15
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = "#"
$ws.Range("K56").Value = 23
$ws.Range("L56").Value = "S4_JPN_WCNPO_OSDWLL_late_Area2"
$ws.Range("L56").Font.ThemeColor = 1

$ws.Range("F57").Value = 15
$ws.Range("G57").Value = 0
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = "#"
$ws.Range("K57").Value = 24
$ws.Range("L57").Value = "S5_TWN_WCNPO_DWLL_late"
$ws.Range("L57").Font.ThemeColor = 1

$ws.Range("F58").Value = 15
$ws.Range("G58").Value = 0
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = "#"
$ws.Range("K58").Value = 25
$ws.Range("L58").Value = "S6_US_WCNPO_LL_deep"
$ws.Range("L58").Font.ThemeColor = 1

$ws.Range("F59").Value = 15
$ws.Range("G59").Value = 0
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = "#"
$ws.Range("K59").Value = 26
$ws.Range("L59").Value = "S7_US_WCNPO_LL_shallow_early"
$ws.Range("L59").Font.ThemeColor = 1

$ws.Range("F60").Value = 15
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = "#"
$ws.Range("K60").Value = 27
$ws.Range("L60").Value = "S8_US_WCNPO_LL_shallow_late"
$ws.Range("L60").Font.ThemeColor = 1

# --- Update the visible selection to span the freshly added block ---
$ws.Range("F34:L60").Select()
